$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D price cells keep their original text formatting
# (values like "38.30" or "52.127.36" would otherwise be auto-converted
# to numbers/dates by Excel, losing trailing zeros or the dotted format).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '52.127.36'
$ws.Range("E2").Value = '  +0.64%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.993.42'
$ws.Range("E3").Value = '  +1.66%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '354.53'
$ws.Range("E5").Value = '  +0.40%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '107.77'
$ws.Range("E6").Value = '  -4.35%  '
$ws.Range("E7").Value = '  +0.18%  '
$ws.Range("E8").Value = '  -0.08%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.624'
$ws.Range("E9").Value = '  -1.38%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '38.30'
$ws.Range("E10").Value = '  -3.13%  '
$ws.Range("E11").Value = '  +1.69%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0857'
$ws.Range("E12").Value = '  -3.67%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.37'
$ws.Range("E13").Value = '  -3.02%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.460.51'
$ws.Range("E14").Value = '  +1.48%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.66'
$ws.Range("E15").Value = '  -2.60%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.994.60'
$ws.Range("E16").Value = '  +2.32%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.02'
$ws.Range("E17").Value = '  +2.99%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '52.190.45'
$ws.Range("E18").Value = '  +0.50%  '
$ws.Range("E19").Value = '  +5.09%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.54'
$ws.Range("E20").Value = '  -1.47%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.68'
$ws.Range("E21").Value = '  -5.06%  '
$ws.Range("E22").Value = '  -1.28%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '69.50'
$ws.Range("E23").Value = '  -2.59%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '263.90'
$ws.Range("E24").Value = '  -2.08%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.180'
$ws.Range("E26").Value = '  +0.53%  '
$ws.Range("E27").Value = '  -0.69%  '
$ws.Range("E28").Value = '  +2.75%  '
$ws.Range("E29").Value = '  -0.07%  '
$ws.Range("E30").Value = '  -1.82%  '
$ws.Range("E31").Value = '  -3.32%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.19'
$ws.Range("E32").Value = '  +0.40%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '36.37'
$ws.Range("E33").Value = '  -3.09%  '
$ws.Range("E34").Value = '  -3.77%  '
$ws.Range("E35").Value = '  -3.85%  '
$ws.Range("E36").Value = '  -1.91%  '
$ws.Range("E37").Value = '  -0.04%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.21'
$ws.Range("E38").Value = '  -3.64%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '17.90'
$ws.Range("E39").Value = '  -4.77%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.97'
$ws.Range("E40").Value = '  -3.68%  '
$ws.Range("E41").Value = '  +1.15%  '
$ws.Range("E42").Value = '  -0.32%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '22.77'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '122.27'
$ws.Range("E44").Value = '  +8.24%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.15'
$ws.Range("E45").Value = '  -2.50%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.125.74'
$ws.Range("E46").Value = '  -1.92%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.38'
$ws.Range("E47").Value = '  -3.98%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.39'
$ws.Range("E48").Value = '  -5.86%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.241'
$ws.Range("E49").Value = '  -1.07%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0333'
$ws.Range("E50").Value = '  -2.71%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.931'
$ws.Range("E51").Value = '  -0.42%  '
